$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New arrivals data for Tuesday, Jan 10 (rows 11-13)
$newRows = @(
    @(10, "Tuesday, Jan 10", "9:00 AM", "FR3693", "Birmingham", "(BHX)", "Ryanair ", "B738", "(EI-DLH)", "8:36 AM", "0 hours, -24 minutes"),
    @(11, "Tuesday, Jan 10", "2:40 PM", "LO3993", "Warsaw", "(WAW)", "LOT ", "E170", "(SP-LDF)", "2:26 PM", "0 hours, -14 minutes"),
    @(12, "Tuesday, Jan 10", "7:30 PM", "W95175", "London", "(LTN)", "Wizz Air ", "A321", "(G-WUKI)", "7:09 PM", "0 hours, -21 minutes")
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 12).Value = $row[10]

    # Columns K and M stay empty but formatted (matching the rest of the table)
    $ws.Cells.Item($r, 11).Borders.LineStyle = 0
    $ws.Cells.Item($r, 13).Borders.LineStyle = 0
}
